$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column J (index 10) is an always-empty column sitting between "breyer" (I)
# and the real "alito" data (K). Delete it so K shifts left into J,
# matching the diff's column removal (K1->J1 "alito", K2:K83 -> J2:J83).
$ws.Columns.Item(10).Delete()

# Many text cells throughout the sheet contain stray leading/trailing
# whitespace (e.g. " joinmajority", "dissent ") that the commit cleans up.
# Trim every text cell in the used range.
$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $v = $cell.Value()
        if ($v -ne $null) {
            if ($v -is [string]) {
                $trimmed = $v.Trim()
                if ($trimmed -ne $v) {
                    $cell.Value = $trimmed
                }
            }
        }
    }
}
